$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Books")

$rows = @(
    @(5, "Ram_Book",    "RAM",     1000),
    @(6, "Java_Book",   "Oracle",  500),
    @(7, "Python_Book", "Python",  500),
    @(8, "Angular_Book","Angular", 900),
    @(9, "React_book",  "React",   5000)
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Range("C10").Select()
